$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Apollo 11")
$ws.Range("E25").Value2 = "LM Launch Mass"
$ws.Range("E30").Value2 = "LM Lunar LO Mass"
$ws.Range("E32").Value2 = "LM Lunar LO RCS"
$ws.Range("E26").Value2 = "LM Descent Empty Calculation"
$ws.Range("E31").Value2 = "LM Ascent Empty Calculation"
$ws.Range("F24").Value2 = "lb"
$ws.Range("E15").Value2 = "Data From Apollo 11 SCOT & Apollo 11 Mission Report"
Write-Host "done"
